$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "(9z,12z)-octadeca-9,12-dienoic acid, deriv."
$ws.Range("B2").Value = "(9z,12z)-octadeca-9,12-dienoic acid"
$ws.Range("C2").Value = "C18H32O2"
$ws.Range("D2").Value = "CCCCCC=CCC=CCCCCCCCC(=O)O"
$ws.Range("E2").Value = 280.4
$ws.Range("F2").Value = 6.8
$ws.Range("G2").Value = "(9z,12z)-octadeca-9,12-dienoic acid"
$ws.Range("H2").Value = 18
$ws.Range("I2").Value = 32
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 0.7710342368045648
$ws.Range("L2").Value = 0.1150356633380885
$ws.Range("M2").Value = 0.1141155492154066
$ws.Range("N2").Value = 17
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0.8396398002853066
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0.1605456490727532
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1.000038998518056

$ws.Range("A3").Value = "9-octadecenoic acid, (z)-, tms derivative"
$ws.Range("B3").Value = "(z)-octadec-9-enoic acid"
$ws.Range("C3").Value = "C18H34O2"
$ws.Range("D3").Value = "CCCCCCCCC=CCCCCCCCC(=O)O"
$ws.Range("E3").Value = 282.5
$ws.Range("F3").Value = 6.5
$ws.Range("G3").Value = "9-octadecenoic acid, (z)-"
$ws.Range("H3").Value = 18
$ws.Range("I3").Value = 34
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 0.7653026548672566
$ws.Range("L3").Value = 0.121316814159292
$ws.Range("M3").Value = 0.1132672566371681
$ws.Range("N3").Value = 17
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.8405345132743363
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0.1593522123893805
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 1.000038998518056

$ws.Range("A4").Value = "4-oxopentanoic acid, deriv."
$ws.Range("B4").Value = "4-oxopentanoic acid"
$ws.Range("C4").Value = "C5H8O3"
$ws.Range("D4").Value = "CC(=O)CCC(=O)O"
$ws.Range("E4").Value = 116.11
$ws.Range("F4").Value = -0.5
$ws.Range("G4").Value = "4-oxopentanoic acid"
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 0.5172250452157436
$ws.Range("L4").Value = 0.06945138230987856
$ws.Range("M4").Value = 0.413375247610025
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 0.1208078546206184
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0.3877099302385669
$ws.Range("W4").Value = 0.491533890276462
$ws.Range("X4").Value = 1.000038998518056

$ws.Range("A5").Value = "benzene-1,2-diol, deriv."
$ws.Range("B5").Value = "benzene-1,2-diol"
$ws.Range("C5").Value = "C6H6O2"
$ws.Range("D5").Value = "C1=CC=C(C(=C1)O)O"
$ws.Range("E5").Value = 110.11
$ws.Range("F5").Value = 0.9
$ws.Range("G5").Value = "benzene-1,2-diol"
$ws.Range("H5").Value = 6
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 0.6544909635818728
$ws.Range("L5").Value = 0.05492689129052766
$ws.Range("M5").Value = 0.2906003087821269
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 6
$ws.Range("P5").Value = 2
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0.6911088911088911
$ws.Range("U5").Value = 0.3089092725456362
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 1.000038998518056

$ws.Range("A6").Value = "9,12-octadecadienoic acid (z,z)-, tms derivative"
$ws.Range("B6").Value = "(9z,12z)-octadeca-9,12-dienoic acid"
$ws.Range("C6").Value = "C18H32O2"
$ws.Range("D6").Value = "CCCCCC=CCC=CCCCCCCCC(=O)O"
$ws.Range("E6").Value = 280.4
$ws.Range("F6").Value = 6.8
$ws.Range("G6").Value = "9,12-octadecadienoic acid (z,z)-"
$ws.Range("H6").Value = 18
$ws.Range("I6").Value = 32
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 0.7710342368045648
$ws.Range("L6").Value = 0.1150356633380885
$ws.Range("M6").Value = 0.1141155492154066
$ws.Range("N6").Value = 17
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0.8396398002853066
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 0.1605456490727532
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 1.000038998518056

$ws.Range("A7").Value = "palmitic acid, tms derivative"
$ws.Range("B7").Value = "hexadecanoic acid"
$ws.Range("C7").Value = "C16H32O2"
$ws.Range("D7").Value = "CCCCCCCCCCCCCCCC(=O)O"
$ws.Range("E7").Value = 256.42
$ws.Range("F7").Value = 6.4
$ws.Range("G7").Value = "palmitic acid"
$ws.Range("H7").Value = 16
$ws.Range("I7").Value = 32
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 0.7494579205990172
$ws.Range("L7").Value = 0.125793619842446
$ws.Range("M7").Value = 0.1247874580765931
$ws.Range("N7").Value = 15
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0.8244793697839481
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0.1755596287341081
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 1.000038998518056

$ws.Range("A8").Value = "phenol, deriv."
$ws.Range("B8").Value = "phenol"
$ws.Range("C8").Value = "C6H6O"
$ws.Range("D8").Value = "C1=CC=C(C=C1)O"
$ws.Range("E8").Value = 94.11
$ws.Range("F8").Value = 1.5
$ws.Range("G8").Value = "phenol"
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 0.765763468281798
$ws.Range("L8").Value = 0.06426522154925088
$ws.Range("M8").Value = 0.1700031877590054
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 6
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0.8193178195728402
$ws.Range("U8").Value = 0.1807140580172139
$ws.Range("V8").Value = 0
$ws.Range("W8").Value = 0
$ws.Range("X8").Value = 1.000038998518056

$ws.Range("A9").Value = "myristic acid, tms derivative"
$ws.Range("B9").Value = "tetradecanoic acid"
$ws.Range("C9").Value = "C14H28O2"
$ws.Range("D9").Value = "CCCCCCCCCCCCCC(=O)O"
$ws.Range("E9").Value = 228.37
$ws.Range("F9").Value = 5.3
$ws.Range("G9").Value = "myristic acid"
$ws.Range("H9").Value = 14
$ws.Range("I9").Value = 28
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 0.7363226343214958
$ws.Range("L9").Value = 0.1235889127293427
$ws.Range("M9").Value = 0.1401147261023777
$ws.Range("N9").Value = 13
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0.8029031834303979
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 0
$ws.Range("V9").Value = 0.1971230897228183
$ws.Range("W9").Value = 0
$ws.Range("X9").Value = 1.000038998518056

$ws.Range("A10").Value = "palmitelaidic acid, tms derivative"
$ws.Range("B10").Value = "(e)-hexadec-9-enoic acid"
$ws.Range("C10").Value = "C16H30O2"
$ws.Range("D10").Value = "CCCCCCC=CCCCCCCCC(=O)O"
$ws.Range("E10").Value = 254.41
$ws.Range("F10").Value = 6.4
$ws.Range("G10").Value = "palmitelaidic acid"
$ws.Range("H10").Value = 16
$ws.Range("I10").Value = 30
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = 0.7553791124562713
$ws.Range("L10").Value = 0.1188632522306513
$ws.Range("M10").Value = 0.1257733579654888
$ws.Range("N10").Value = 15
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0.8230690617507173
$ws.Range("T10").Value = 0
$ws.Range("U10").Value = 0
$ws.Range("V10").Value = 0.1769466609016941
$ws.Range("W10").Value = 0
$ws.Range("X10").Value = 1.000038998518056

$ws.Range("A11").Value = "9-octadecenoic acid, (e)-, deriv."
$ws.Range("B11").Value = "(e)-octadec-9-enoic acid"
$ws.Range("C11").Value = "C18H34O2"
$ws.Range("D11").Value = "CCCCCCCCC=CCCCCCCCC(=O)O"
$ws.Range("E11").Value = 282.5
$ws.Range("F11").Value = 6.5
$ws.Range("G11").Value = "9-octadecenoic acid, (e)-"
$ws.Range("H11").Value = 18
$ws.Range("I11").Value = 34
$ws.Range("J11").Value = 2
$ws.Range("K11").Value = 0.7653026548672566
$ws.Range("L11").Value = 0.121316814159292
$ws.Range("M11").Value = 0.1132672566371681
$ws.Range("N11").Value = 17
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0.8405345132743363
$ws.Range("T11").Value = 0
$ws.Range("U11").Value = 0
$ws.Range("V11").Value = 0.1593522123893805
$ws.Range("W11").Value = 0
$ws.Range("X11").Value = 1.000038998518056

$ws.Range("A12").Value = "benzoic acid, deriv."
$ws.Range("B12").Value = "benzoic acid"
$ws.Range("C12").Value = "C7H6O2"
$ws.Range("D12").Value = "C1=CC=C(C=C1)C(=O)O"
$ws.Range("E12").Value = 122.12
$ws.Range("F12").Value = 1.9
$ws.Range("G12").Value = "benzoic acid"
$ws.Range("H12").Value = 7
$ws.Range("I12").Value = 6
$ws.Range("J12").Value = 2
$ws.Range("K12").Value = 0.6884785456927611
$ws.Range("L12").Value = 0.04952505732066819
$ws.Range("M12").Value = 0.2620209629872257
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 6
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("T12").Value = 0.6313953488372093
$ws.Range("U12").Value = 0
$ws.Range("V12").Value = 0.3686292171634458
$ws.Range("W12").Value = 0
$ws.Range("X12").Value = 1.000038998518056

$ws.Range("A13").Value = "hexadecanoic acid, deriv."
$ws.Range("B13").Value = "hexadecanoic acid"
$ws.Range("C13").Value = "C16H32O2"
$ws.Range("D13").Value = "CCCCCCCCCCCCCCCC(=O)O"
$ws.Range("E13").Value = 256.42
$ws.Range("F13").Value = 6.4
$ws.Range("G13").Value = "hexadecanoic acid"
$ws.Range("H13").Value = 16
$ws.Range("I13").Value = 32
$ws.Range("J13").Value = 2
$ws.Range("K13").Value = 0.7494579205990172
$ws.Range("L13").Value = 0.125793619842446
$ws.Range("M13").Value = 0.1247874580765931
$ws.Range("N13").Value = 15
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0.8244793697839481
$ws.Range("T13").Value = 0
$ws.Range("U13").Value = 0
$ws.Range("V13").Value = 0.1755596287341081
$ws.Range("W13").Value = 0
$ws.Range("X13").Value = 1.000038998518056
